$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values in column F (dSF) - repulled data
$ws.Range("F2").Value = -8
$ws.Range("F4").Value = 0
$ws.Range("F5").Value = 4
$ws.Range("F10").Value = 5
$ws.Range("F13").Value = -2
$ws.Range("F16").Value = -4
$ws.Range("F17").Value = 0
$ws.Range("F23").Value = -2
$ws.Range("F25").Value = -7
$ws.Range("F29").Value = 0
$ws.Range("F30").Value = 3
$ws.Range("F34").Value = -3
$ws.Range("F35").Value = 1
$ws.Range("F41").Value = -2
$ws.Range("F42").Value = -4
$ws.Range("F43").Value = -3
$ws.Range("F57").Value = -1
$ws.Range("F60").Value = -4
$ws.Range("F62").Value = -7
